# Slide 11 ("Contributions"), shape "TextBox 5" (3rd shape) holds the
# contributor list.  Two edits land in this same text frame:
#
#   1. The "Tester: ..." paragraph was re-typed as a single run instead
#      of four separate runs (the split runs collapse into one string).
#   2. The "Project Manager: Aaditya S Shah" paragraph had a leading
#      space removed and now starts with its own "Project " run,
#      followed by a "Manager: Aaditya S Shah " run.

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(11)
$shp = $s.Shapes.Item(3)
$tr  = $shp.TextFrame.TextRange

# The shape autosizes (spAutoFit) whenever its text changes; remember the
# original height so it can be restored once both edits are done.
$origHeight = $shp.Height

# --- 1. Merge the "Tester: ..." runs into a single run ---------------
$fullText      = $tr.Text
$testerStart0  = $fullText.IndexOf("Tester: Krishna S")
$testerEnd0    = $fullText.IndexOf("Hariharan.") + "Hariharan.".Length
$testerStart1  = $testerStart0 + 1
$testerLen     = $testerEnd0 - $testerStart0

$testerRange = $tr.Characters($testerStart1, $testerLen)
$testerRange.Text = "Tester: Krishna S, Amitesh M, Selva Akash M , Abhijit S , Hariharan."

# --- 2. Drop the leading space and split "Project " into its own run -
$fullText2   = $tr.Text
$pmStart0    = $fullText2.IndexOf(" Project Manager: Aaditya S Shah ")
$pmStart1    = $pmStart0 + 1

$leadingSpace = $tr.Characters($pmStart1, 1)
$leadingSpace.Text = ""

$projectRun = $tr.Characters($pmStart1, "Project ".Length)
$projectRun.Font.Bold = $true
$projectRun.Font.Size = 36

# Restore the shape's original (pre-autofit-recalc) height.
$shp.Height = $origHeight
